# Generate Report for Handback
#
# The localization-status workbook tracks handoff/handback round trips for
# each localized file. Both the "zh-cn" and "de-de" worksheets list the two
# tracked source files (rows 2 and 3) and a sentinel ".localization-config"
# row (row 4). This handback event marks both tracked files as handed back
# "in sync with en-US", and records the target/handback file names plus the
# handback timestamp for each locale.

$wb = $excel.ActiveWorkbook

function Apply-Handback {
    param(
        [object]$ws,
        [int]$row,
        [string]$mdDisplay,
        [string]$mdUrl,
        [string]$xlfDisplay,
        [string]$xlfUrl,
        [string]$handbackTime
    )

    # Status column: the file is no longer just "Ready for handoff" - it has
    # come back from the loc vendor in sync with the English source.
    $ws.Range("B$row").Value = "Handed back: in sync with en-US"

    # "Latest Target File" (E) and "Latest Handback File" (F) now point at
    # the same source markdown / xlf pair recorded under "Latest Handoff
    # File" (columns A/C), carried over with their original hyperlinks.
    $ws.Hyperlinks.Add($ws.Range("E$row"), $mdUrl, "", "", $mdDisplay)
    $ws.Hyperlinks.Add($ws.Range("F$row"), $xlfUrl, "", "", $xlfDisplay)

    # "Latest Handback DateTime" (G) records when the handback happened,
    # replacing the 0001-01-01 placeholder.
    $ws.Range("G$row").Value = $handbackTime
}

$wsZh = $wb.Worksheets.Item("zh-cn")

Apply-Handback $wsZh 2 `
    "de9964be-e953-48bb-b279-f24b6c27c2a8.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a97d1dfcda1b7eb8d1b3118fff59927266c48d9e/e2e/de9964be-e953-48bb-b279-f24b6c27c2a8.md" `
    "de9964be-e953-48bb-b279-f24b6c27c2a8.02ab420233d2233a0290be36e2cb31c875337c9d.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/774d3748ef0fd2ad5ed0123753c68aa84f285a28/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/de9964be-e953-48bb-b279-f24b6c27c2a8.02ab420233d2233a0290be36e2cb31c875337c9d.zh-cn.xlf" `
    "2016-03-09 16:11:45"

Apply-Handback $wsZh 3 `
    "f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a97d1dfcda1b7eb8d1b3118fff59927266c48d9e/e2e/f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md" `
    "f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.6ada2efc41a576ca0a5ca21d912e7e1a1b1a8af4.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/774d3748ef0fd2ad5ed0123753c68aa84f285a28/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.6ada2efc41a576ca0a5ca21d912e7e1a1b1a8af4.zh-cn.xlf" `
    "2016-03-09 16:11:45"

$wsDe = $wb.Worksheets.Item("de-de")

Apply-Handback $wsDe 2 `
    "de9964be-e953-48bb-b279-f24b6c27c2a8.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a97d1dfcda1b7eb8d1b3118fff59927266c48d9e/e2e/de9964be-e953-48bb-b279-f24b6c27c2a8.md" `
    "de9964be-e953-48bb-b279-f24b6c27c2a8.02ab420233d2233a0290be36e2cb31c875337c9d.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f17071bd68f81d30e54dd6a2b0f0ad5a113cb07f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/de9964be-e953-48bb-b279-f24b6c27c2a8.02ab420233d2233a0290be36e2cb31c875337c9d.de-de.xlf" `
    "2016-03-09 16:12:07"

Apply-Handback $wsDe 3 `
    "f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md" `
    "https://github.com/OpenLocalizationTest/oltest/blob/a97d1dfcda1b7eb8d1b3118fff59927266c48d9e/e2e/f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.md" `
    "f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.6ada2efc41a576ca0a5ca21d912e7e1a1b1a8af4.de-de.xlf" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f17071bd68f81d30e54dd6a2b0f0ad5a113cb07f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/f03ca9ea-e27c-4d6c-87fb-ef8cd8f0cdbe.6ada2efc41a576ca0a5ca21d912e7e1a1b1a8af4.de-de.xlf" `
    "2016-03-09 16:12:07"
